$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @(
    45677.5626052291,
    45524.56260522935,
    45557.56260522941,
    45460.56260522945,
    45565.56260522948,
    45504.56260522951,
    45507.56260522956,
    45642.56260522959,
    45448.56260522969,
    45665.56260522972,
    45776.56260522976,
    45646.56260522978,
    45713.5626052298,
    45804.56260522984,
    45692.56260522986,
    45514.56260522988,
    45632.56260522992,
    45703.56260522994,
    45588.56260522998,
    45656.56260523,
    45702.56260523002,
    45786.56260523006,
    45581.56260523009,
    45488.56260523012,
    45603.56260523015,
    45812.56260523017,
    45586.5626052302,
    45653.56260523023,
    45503.56260523025,
    45694.56260523029,
    45538.56260523031,
    45700.56260523035,
    45575.56260523037,
    45785.56260523041,
    45484.56260523044,
    45811.56260523047,
    45694.56260523049,
    45677.56260523052,
    45731.56260523055,
    45571.56260523057,
    45479.56260523059,
    45457.56260523063,
    45449.56260523065,
    45466.56260523068,
    45807.56260523071,
    45566.56260523074,
    45478.56260523076,
    45649.5626052308,
    45686.56260523082,
    45775.56260523084,
    45612.56260523088,
    45667.5626052309,
    45486.56260523092,
    45813.56260523096,
    45746.56260523098,
    45527.562605231,
    45667.56260523104,
    45592.56260523106,
    45610.5626052311,
    45737.56260523112,
    45523.56260523115,
    45686.56260523118,
    45621.56260523121,
    45594.56260523123,
    45606.56260523126,
    45725.56260523129,
    45602.56260523132,
    45687.56260523135,
    45769.56260523137,
    45752.5626052314,
    45613.56260523143,
    45597.56260523145,
    45687.56260523148,
    45687.56260523151,
    45578.56260523153,
    45513.56260523155,
    45473.56260523159,
    45450.56260523161,
    45639.56260523164,
    45565.56260523167,
    45735.5626052317,
    45699.56260523175,
    45710.5626052318,
    45798.56260523184,
    45671.56260523188,
    45674.56260523193,
    45788.56260523196,
    45774.56260523201,
    45760.56260523206,
    45458.56260523212,
    45723.56260523215,
    45504.5626052322,
    45627.56260523226,
    45795.56260523231,
    45458.56260523234,
    45803.56260523237,
    45726.56260523239,
    45712.56260523242,
    45489.56260523245,
    45748.56260523247
)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $newDates[$i]
}
